$d = $word.ActiveDocument

$replacements = @(
    @("836×2=", "759×9="),
    @("913×9=", "827×9="),
    @("952×4=", "505×9="),
    @("751×4=", "774×5="),
    @("900×3=", "120×8="),
    @("150×7=", "692×9="),
    @("701×7=", "306×4="),
    @("113×6=", "611×4="),
    @("399×6=", "240×4="),
    @("588×5=", "870×7="),
    @("160×2=", "323×4="),
    @("673×4=", "596×6="),
    @("424×7=", "656×6="),
    @("239×5=", "529×9="),
    @("386×4=", "947×3="),
    @("910×2=", "845×9="),
    @("580×6=", "331×3="),
    @("306×2=", "354×3="),
    @("747×4=", "975×2="),
    @("748×8=", "104×2="),
    @("872×7=", "793×8="),
    @("744×9=", "510×5="),
    @("310×3=", "778×5="),
    @("674×9=", "423×7="),
    @("350×8=", "638×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
